$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# 2. Remove the old table (rows 5-10, columns A:D) completely so the rows
#    disappear from the sheet rather than leaving empty cells behind.
$ws.Range("A5:D10").Clear()

# 3. New bold+underlined "Source Type" label above the table.
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# 4. Header row with column labels.
$ws.Range("B9").Value = "Micro"
$ws.Range("B9").Style = "title"
$ws.Range("C9").Value = "SMEs"
$ws.Range("C9").Style = "title"
$ws.Range("D9").Value = "MSMEs"
$ws.Range("D9").Style = "title"

# 5. New "Employment (% of total)" data row.
$ws.Range("A10").Value = "Employment (% of total)"
$ws.Range("A10").Style = "title"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.6"
$ws.Range("D10").Style = "Normal"

# 6. Enterprises (absolute #) row.
$ws.Range("A11").Value = "Enterprises (absolute #)"
$ws.Range("A11").Style = "title"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "217250"
$ws.Range("D11").Style = "Normal"

# 7. Enterprises density row.
$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("A12").Style = "title"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.6"
$ws.Range("D12").Style = "Normal"

# 8. Employment (absolute #) row.
$ws.Range("A13").Value = "Employment (absolute #)"
$ws.Range("A13").Style = "title"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "492181"
$ws.Range("D13").Style = "Normal"

# 9. Enterprises (% of total) row.
$ws.Range("A14").Value = "Enterprises (% of total)"
$ws.Range("A14").Style = "title"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.9"
$ws.Range("D14").Style = "Normal"

# 10. Source line (italic).
$ws.Range("A15").Value = "Source: DGEEC, 2010"
$ws.Range("A15").Style = "source"

# 11. Additional attribution block further down the sheet.
$ws.Range("A23").Value = "DGEEC"
$ws.Range("A23").Style = "title"

$ws.Range("A24").Value = 'Dirección General de Estadística, Encuestas y Censos (DGEEC), "Censo Económico Nacional 2011", 2013, p.57. Available at http://www.dgeec.gov.py/Publicaciones/Biblioteca/CEN2011/resultados_finales_CEN.pdf'
$ws.Range("A24").Style = "source"
